$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.779.29"
$ws.Range("E2").Value = "  +0.63%  "
$ws.Range("D3").Value = "'3.164.70"
$ws.Range("E3").Value = "  +0.69%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'616.37"
$ws.Range("E5").Value = "  +2.26%  "
$ws.Range("D6").Value = "'147.17"
$ws.Range("E6").Value = "  -1.88%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").Value = "'3.162.40"
$ws.Range("E8").Value = "  +0.68%  "
$ws.Range("E9").Value = "  -1.13%  "
$ws.Range("E10").Value = "  -1.11%  "
$ws.Range("E11").Value = "  -2.08%  "
$ws.Range("E12").Value = "  -1.31%  "
$ws.Range("E13").Value = "  -0.42%  "
$ws.Range("D14").Value = "'35.84"
$ws.Range("E14").Value = "  -3.37%  "
$ws.Range("D15").Value = "'3.686.06"
$ws.Range("E15").Value = "  +1.91%  "
$ws.Range("E16").Value = "  +2.89%  "
$ws.Range("D17").Value = "'64.752.85"
$ws.Range("E17").Value = "  +0.45%  "
$ws.Range("D18").Value = "'3.161.26"
$ws.Range("E18").Value = "  +0.75%  "
$ws.Range("E19").Value = "  -1.90%  "
$ws.Range("D20").Value = "'479.24"
$ws.Range("E20").Value = "  -1.03%  "
$ws.Range("D21").Value = "'14.74"
$ws.Range("E21").Value = "  -0.37%  "
$ws.Range("E22").Value = "  +1.21%  "
$ws.Range("D23").Value = "'7.94"
$ws.Range("E23").Value = "  +1.32%  "
$ws.Range("D24").Value = "'13.80"
$ws.Range("E24").Value = "  -1.20%  "
$ws.Range("D25").Value = "'84.48"
$ws.Range("E25").Value = "  -0.58%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("E27").Value = "  -3.63%  "
$ws.Range("D28").Value = "'8.58"
$ws.Range("E28").Value = "  -1.61%  "
$ws.Range("D29").Value = "'6.93"
$ws.Range("E29").Value = "  -3.40%  "
$ws.Range("E30").Value = "  -7.26%  "
$ws.Range("D31").Value = "'2.08"
$ws.Range("E31").Value = "  -8.10%  "
$ws.Range("E32").Value = "  +0.23%  "
$ws.Range("D33").Value = "'2.70"
$ws.Range("E33").Value = "  -0.92%  "
$ws.Range("E34").Value = "  -1.16%  "
$ws.Range("D35").Value = "'1.14"
$ws.Range("E35").Value = "  +2.54%  "
$ws.Range("D36").Value = "0.0₃0779"
$ws.Range("E36").Value = "  +2.28%  "
$ws.Range("E37").Value = "  -2.05%  "
$ws.Range("B38").Value = "OKB"
$ws.Range("C38").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D38").Value = "'53.04"
$ws.Range("E38").Value = "  -3.10%  "
$ws.Range("B39").Value = "dogwifhat"
$ws.Range("C39").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D39").Value = "'3.17"
$ws.Range("E39").Value = "  -5.01%  "
$ws.Range("D40").Value = "'460.24"
$ws.Range("E40").Value = "  +1.57%  "
$ws.Range("E41").Value = "  -1.27%  "
$ws.Range("D42").Value = "'0.120"
$ws.Range("E42").Value = "  -4.14%  "
$ws.Range("D43").Value = "'8.41"
$ws.Range("E43").Value = "  -1.82%  "
$ws.Range("D44").Value = "'2.849.25"
$ws.Range("E44").Value = "  -1.82%  "
$ws.Range("D45").Value = "'2.33"
$ws.Range("E45").Value = "  -4.55%  "
$ws.Range("E46").Value = "  -2.79%  "
$ws.Range("E47").Value = "  +3.61%  "
$ws.Range("D48").Value = "'26.60"
$ws.Range("E48").Value = "  -1.87%  "
$ws.Range("E49").Value = "  +0.12%  "
$ws.Range("E50").Value = "  -1.79%  "
$ws.Range("D51").Value = "'120.32"
$ws.Range("E51").Value = "  +0.59%  "

Write-Output "Applied updates to cryptos list."
